# Add a new column W ("param_E_pv3_solar") to Sheet1, and update columns
# U (param_P_to_charging_station1) and V (param_P_to_charging_station2)
# with the new simulation results, per the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) -------------------------------------------------
# Clone the formatting of V1 (bold, centered, bordered header style) onto
# the new W1 cell, then give it its own label.
$ws.Range("V1").Copy($ws.Range("W1"))
$ws.Range("W1").Value = "param_E_pv3_solar"

# --- Data rows (rows 2-17) ----------------------------------------------
# Column U: every row becomes 0.12
# Column V: updated simulation values (several rows change)
# Column W: brand-new column of simulation values

$uValues = @{
    2  = 0.12
    3  = 0.12
    4  = 0.12
    5  = 0.12
    6  = 0.12
    7  = 0.12
    8  = 0.12
    9  = 0.12
    10 = 0.12
    11 = 0.12
    12 = 0.12
    13 = 0.12
    14 = 0.12
    15 = 0.12
    16 = 0.12
    17 = 0.12
}

$vValues = @{
    2  = 0
    3  = 87.94500000000002
    4  = 161.8688194444445
    5  = 0
    6  = 0
    7  = 0
    8  = 59.89148611111113
    9  = 54.11648611111113
    10 = 0
    11 = 71.51084722222224
    12 = 30.25000000000002
    13 = 0
    14 = 0
    15 = 0
    16 = 0
    17 = 0
}

$wValues = @{
    2  = 0
    3  = 44.99000000000002
    4  = 113.5164861111111
    5  = 0
    6  = 47.85000000000002
    7  = 109.395
    8  = 49.33500000000002
    9  = 0
    10 = 0
    11 = 0
    12 = 0
    13 = 0
    14 = 0
    15 = 0
    16 = 0
    17 = 0
}

for ($row = 2; $row -le 17; $row++) {
    $ws.Cells.Item($row, 21).Value = $uValues[$row]   # column U
    $ws.Cells.Item($row, 22).Value = $vValues[$row]   # column V
    $ws.Cells.Item($row, 23).Value = $wValues[$row]   # column W
}
